# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Palta, Hass - "1a nueva(o)" / "2a nueva(o)")
# right before the existing row 224, pushing the rest of the table down by two
# rows (old A1:T313 -> new A1:T315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 224:225 (existing rows 224+ shift down to 226+).
$ws.Rows("224:225").Insert()

$d = Get-Date -Year 2021 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0

# New row 224 — Palta, Hass, "1a nueva(o)"
$ws.Cells.Item(224,1).Value = 7
$ws.Cells.Item(224,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(224,3).Value = "Ñuble"
$ws.Cells.Item(224,4).Value = $d
$ws.Cells.Item(224,5).Value = 16
$ws.Cells.Item(224,6).Value = "Fruta"
$ws.Cells.Item(224,7).Value = 100106
$ws.Cells.Item(224,8).Value = "Oleaginosos"
$ws.Cells.Item(224,9).Value = 100106002
$ws.Cells.Item(224,10).Value = "Palta"
$ws.Cells.Item(224,11).Value = "Hass"
$ws.Cells.Item(224,12).Value = "1a nueva(o)"
$ws.Cells.Item(224,13).Value = 60
$ws.Cells.Item(224,14).Value = 2800
$ws.Cells.Item(224,15).Value = 2800
$ws.Cells.Item(224,16).Value = 2800
$ws.Cells.Item(224,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(224,18).Value = "Provincia de Quillota"
$ws.Cells.Item(224,19).Value = 2800
$ws.Cells.Item(224,20).Value = 1

# New row 225 — Palta, Hass, "2a nueva(o)"
$ws.Cells.Item(225,1).Value = 7
$ws.Cells.Item(225,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(225,3).Value = "Ñuble"
$ws.Cells.Item(225,4).Value = $d
$ws.Cells.Item(225,5).Value = 16
$ws.Cells.Item(225,6).Value = "Fruta"
$ws.Cells.Item(225,7).Value = 100106
$ws.Cells.Item(225,8).Value = "Oleaginosos"
$ws.Cells.Item(225,9).Value = 100106002
$ws.Cells.Item(225,10).Value = "Palta"
$ws.Cells.Item(225,11).Value = "Hass"
$ws.Cells.Item(225,12).Value = "2a nueva(o)"
$ws.Cells.Item(225,13).Value = 60
$ws.Cells.Item(225,14).Value = 2600
$ws.Cells.Item(225,15).Value = 2600
$ws.Cells.Item(225,16).Value = 2600
$ws.Cells.Item(225,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(225,18).Value = "Provincia de Quillota"
$ws.Cells.Item(225,19).Value = 2600
$ws.Cells.Item(225,20).Value = 1

Write-Output "done"
